# Reorder data rows 2-21 (the full data block, columns A:AY) of the
# "Artfynd" worksheet according to the row permutation observed between
# the pre- and post-edit workbook snapshots. Only the row order changes;
# no individual cell value is altered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 21
$lastCol = "AY"

# Columns Y and AA hold date-like text (e.g. "2021-08-11") stored as plain
# text in the workbook. Force them to stay text so that re-assigning the
# values through Value2 does not let Excel auto-convert them to date
# serial numbers.
$ws.Range("Y" + $firstRow + ":Y" + $lastRow).NumberFormat = "@"
$ws.Range("AA" + $firstRow + ":AA" + $lastRow).NumberFormat = "@"

$srcRange = $ws.Range("A" + $firstRow + ":" + $lastCol + $lastRow)
$data = $srcRange.Value2

$numRows = $data.GetUpperBound(0)
$numCols = $data.GetUpperBound(1)

# Mapping of new (1-based, relative) row index -> old (1-based, relative) row index
$rowMap = @{
    1  = 11
    2  = 1
    3  = 12
    4  = 13
    5  = 14
    6  = 15
    7  = 16
    8  = 2
    9  = 3
    10 = 17
    11 = 18
    12 = 4
    13 = 5
    14 = 19
    15 = 20
    16 = 6
    17 = 7
    18 = 8
    19 = 9
    20 = 10
}

$newData = New-Object 'object[,]' $numRows, $numCols

for ($newIdx = 1; $newIdx -le $numRows; $newIdx++) {
    $oldIdx = $rowMap[$newIdx]
    for ($c = 1; $c -le $numCols; $c++) {
        $newData[$newIdx - 1, $c - 1] = $data[$oldIdx, $c]
    }
}

$srcRange.Value2 = $newData
